{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies the same logical change as the target XML diff:\n//  1. Insert a new \"Mateusz Kosiba - Sublime Text 3\" paragraph (plus a\n//     trailing blank paragraph) right after the \"Jarek Porada - Sublime\n//     Text 3\" paragraph.\n//  2. Drop the stale <w:lastRenderedPageBreak/> marker that sits in front\n//     of the \"Baza danych:\" run (Word regenerates/removes these on edit).\n//  3. Rewrite the final \"Konwencja nazewnicza\" paragraph (fixing the\n//     \"Konwe\" + \"ncja...\" split into a single run) and append the new\n//     \"Cel aplikacji\" section content that follows it, moving the\n//     _GoBack bookmark down to the new final run.\n\nconst NS_PKG =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\\n' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>';\nconst NS_PKG_END = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>';\n\nfunction wrap(innerXml) {\n  return NS_PKG + innerXml + NS_PKG_END;\n}\n\nconst body = context.document.body;\n\n// ---------------------------------------------------------------------\n// Step 1: insert \"Mateusz Kosiba\" paragraph (+ blank paragraph) after the\n// \"Jarek Porada - Sublime Text 3\" paragraph.\n// ---------------------------------------------------------------------\nconst jarekResults = body.search(\"Jarek Porada\", { matchCase: true });\njarekResults.load(\"items\");\nawait context.sync();\n\nif (jarekResults.items.length === 0) {\n  throw new Error(\"Could not find 'Jarek Porada' paragraph anchor\");\n}\n\nconst jarekParagraph = jarekResults.items[0].paragraphs.getFirst();\nconst jarekRange = jarekParagraph.getRange(\"Whole\");\n\nconst mateuszXml = wrap(\n  '<w:p><w:pPr><w:ind w:firstLine=\"720\"/><w:rPr><w:lang w:val=\"pl-PL\"/></w:rPr></w:pPr>' +\n  '<w:r><w:rPr><w:lang w:val=\"pl-PL\"/></w:rPr><w:t xml:space=\"preserve\">Mateusz Kosiba \\u2013 </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:rPr><w:lang w:val=\"pl-PL\"/></w:rPr><w:t>Sublime</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:rPr><w:lang w:val=\"pl-PL\"/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:rPr><w:lang w:val=\"pl-PL\"/></w:rPr><w:t>Text</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:rPr><w:lang w:val=\"pl-PL\"/></w:rPr><w:t xml:space=\"preserve\"> 3 </w:t></w:r>' +\n  '</w:p>' +\n  '<w:p><w:pPr><w:ind w:firstLine=\"720\"/><w:rPr><w:lang w:val=\"pl-PL\"/></w:rPr></w:pPr></w:p>'\n);\n\njarekRange.insertOoxml(mateuszXml, Word.InsertLocation.after);\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// Step 2: remove the stale lastRenderedPageBreak in front of \"Baza\n// danych:\" by rewriting that run from scratch (keeping its rPr).\n// ---------------------------------------------------------------------\nconst bazaResults = body.search(\"Baza danych:\", { matchCase: true });\nbazaResults.load(\"items\");\nawait context.sync();\n\nif (bazaResults.items.length === 0) {\n  throw new Error(\"Could not find 'Baza danych:' run anchor\");\n}\n\nconst bazaXml = wrap(\n  '<w:p><w:r><w:rPr><w:lang w:val=\"pl-PL\"/></w:rPr><w:t>Baza danych:</w:t></w:r></w:p>'\n);\nbazaResults.items[0].insertOoxml(bazaXml, Word.InsertLocation.replace);\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// Step 3: rewrite the trailing \"Konwencja nazewnicza\" paragraph (merging\n// the split \"Konwe\"/\"ncja...\" runs and dropping the bookmark from here)\n// and append the new \"Cel aplikacji\" paragraphs after it, with the\n// _GoBack bookmark re-attached to the new final run.\n// ---------------------------------------------------------------------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nconst lastText = lastParagraph.text || \"\";\nif (lastText.indexOf(\"Konwe\") === -1) {\n  throw new Error(\"Could not find trailing 'Konwencja nazewnicza' paragraph\");\n}\nconst lastRange = lastParagraph.getRange(\"Whole\");\n\nconst tailXml = wrap(\n  '<w:p><w:pPr><w:rPr><w:lang w:val=\"pl-PL\"/></w:rPr></w:pPr>' +\n  '<w:r><w:rPr><w:lang w:val=\"pl-PL\"/></w:rPr><w:t xml:space=\"preserve\">Konwencja nazewnicza zgodna ze standardami pisania kodu w </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:rPr><w:lang w:val=\"pl-PL\"/></w:rPr><w:t>Pythonie</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:rPr><w:lang w:val=\"pl-PL\"/></w:rPr><w:t xml:space=\"preserve\"> (PEP8).</w:t></w:r>' +\n  '</w:p>' +\n\n  '<w:p><w:pPr><w:rPr><w:lang w:val=\"pl-PL\"/></w:rPr></w:pPr>' +\n  '<w:r><w:rPr><w:lang w:val=\"pl-PL\"/></w:rPr><w:t xml:space=\"preserve\">Nie definiujemy konwencji  programowania dla bazy danych poniewa\u017c dla generowanych modeli w </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:rPr><w:lang w:val=\"pl-PL\"/></w:rPr><w:t>Django</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:rPr><w:lang w:val=\"pl-PL\"/></w:rPr><w:t xml:space=\"preserve\"> nie istnieje potrzeba pisania w j\u0119zyku SQL.</w:t></w:r>' +\n  '</w:p>' +\n\n  '<w:p><w:pPr><w:rPr><w:lang w:val=\"pl-PL\"/></w:rPr></w:pPr></w:p>' +\n\n  '<w:p><w:pPr><w:rPr><w:sz w:val=\"28\"/><w:lang w:val=\"pl-PL\"/></w:rPr></w:pPr>' +\n  '<w:r><w:rPr><w:sz w:val=\"28\"/><w:lang w:val=\"pl-PL\"/></w:rPr><w:t>5. Cel aplikacji</w:t></w:r>' +\n  '</w:p>' +\n\n  '<w:p><w:pPr><w:rPr><w:sz w:val=\"28\"/><w:lang w:val=\"pl-PL\"/></w:rPr></w:pPr></w:p>' +\n\n  '<w:p><w:pPr><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/><w:lang w:val=\"pl-PL\"/></w:rPr></w:pPr>' +\n  '<w:r><w:rPr><w:sz w:val=\"24\"/><w:lang w:val=\"pl-PL\"/></w:rPr>' +\n  '<w:t xml:space=\"preserve\">Aplikacja ma na celu informowa\u0107 potencjalnych zawodnik\u00f3w o mo\u017cliwo\u015bci przyst\u0105pienia do zawod\u00f3w. Wspomaga\u0107 rejestracj\u0119 uczestnik\u00f3w. Informowa\u0107 na bie\u017c\u0105co o tym jakie s\u0105 wyniki </w:t></w:r>' +\n  '<w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/><w:lang w:val=\"pl-PL\"/></w:rPr>' +\n  '<w:t xml:space=\"preserve\">zawod\u00f3w. Umo\u017cliwia\u0107 agregacj\u0119 tre\u015bci zwi\u0105zanych z zawodami. </w:t></w:r>' +\n  '</w:p>' +\n\n  '<w:p><w:pPr><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/><w:lang w:val=\"pl-PL\"/></w:rPr></w:pPr></w:p>' +\n\n  '<w:p><w:pPr><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/><w:lang w:val=\"pl-PL\"/></w:rPr></w:pPr></w:p>' +\n\n  '<w:p><w:pPr><w:rPr><w:sz w:val=\"28\"/><w:szCs w:val=\"24\"/><w:lang w:val=\"pl-PL\"/></w:rPr></w:pPr>' +\n  '<w:r><w:rPr><w:sz w:val=\"28\"/><w:szCs w:val=\"24\"/><w:lang w:val=\"pl-PL\"/></w:rPr><w:t xml:space=\"preserve\">6. </w:t></w:r>' +\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>' +\n  '</w:p>' +\n\n  '<w:p><w:pPr><w:rPr><w:lang w:val=\"pl-PL\"/></w:rPr></w:pPr></w:p>'\n);\n\nlastRange.insertOoxml(tailXml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Applies the same logical change as the target XML diff:\n#  1. Insert a new \"Mateusz Kosiba - Sublime Text 3\" paragraph (plus a\n#     trailing blank paragraph) right after the \"Jarek Porada - Sublime\n#     Text 3\" paragraph.\n#  2. Drop the stale <w:lastRenderedPageBreak/> marker that sits in front\n#     of the \"Baza danych:\" run (Word regenerates/removes these on edit).\n#  3. Rewrite the final \"Konwencja nazewnicza\" paragraph (fixing the\n#     \"Konwe\" + \"ncja...\" split into a single run) and append the new\n#     \"Cel aplikacji\" section content that follows it, moving the\n#     _GoBack bookmark down to the new final run.\n\n$d = $word.ActiveDocument\n\nfunction New-PkgXml([string]$bodyXml) {\n  return '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' + $bodyXml + '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>'\n}\n\n# ---------------------------------------------------------------------\n# Step 1: insert \"Mateusz Kosiba\" paragraph (+ blank paragraph) after the\n# \"Jarek Porada - Sublime Text 3\" paragraph.\n# ---------------------------------------------------------------------\n$findRng = $d.Content\n$find = $findRng.Find\n$find.Text = \"Jarek Porada\"\n$find.MatchCase = $true\n$found = $find.Execute()\nif (-not $found) {\n  throw \"Could not find 'Jarek Porada' paragraph anchor\"\n}\n$jarekPara = $findRng.Paragraphs(1)\n$jarekRange = $jarekPara.Range\n\n$mateuszXml = New-PkgXml (\n  '<w:p><w:pPr><w:ind w:firstLine=\"720\"/><w:rPr><w:lang w:val=\"pl-PL\"/></w:rPr></w:pPr>' +\n  '<w:r><w:rPr><w:lang w:val=\"pl-PL\"/></w:rPr><w:t xml:space=\"preserve\">Mateusz Kosiba \\u2013 </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:rPr><w:lang w:val=\"pl-PL\"/></w:rPr><w:t>Sublime</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:rPr><w:lang w:val=\"pl-PL\"/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:rPr><w:lang w:val=\"pl-PL\"/></w:rPr><w:t>Text</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:rPr><w:lang w:val=\"pl-PL\"/></w:rPr><w:t xml:space=\"preserve\"> 3 </w:t></w:r>' +\n  '</w:p>' +\n  '<w:p><w:pPr><w:ind w:firstLine=\"720\"/><w:rPr><w:lang w:val=\"pl-PL\"/></w:rPr></w:pPr></w:p>'\n) -replace '\\\\u2013', [char]0x2013\n\n# Create a fresh empty paragraph right after \"Jarek Porada...\" and replace\n# its content with the two new paragraphs above (keeps paragraph marks\n# intact instead of merging runs into the neighbouring paragraph).\n$jarekRange.InsertParagraphAfter()\n$placeholder = $findRng.Paragraphs(1).Next()\n$placeholder.Range.InsertXML($mateuszXml)\n\n# ---------------------------------------------------------------------\n# Step 2: remove the stale lastRenderedPageBreak in front of \"Baza\n# danych:\" by rewriting that paragraph from scratch (keeping its rPr).\n# ---------------------------------------------------------------------\n$bazaFindRng = $d.Content\n$bazaFind = $bazaFindRng.Find\n$bazaFind.Text = \"Baza danych:\"\n$bazaFind.MatchCase = $true\n$bazaFound = $bazaFind.Execute()\nif (-not $bazaFound) {\n  throw \"Could not find 'Baza danych:' paragraph anchor\"\n}\n$bazaPara = $bazaFindRng.Paragraphs(1)\n$bazaXml = New-PkgXml (\n  '<w:p><w:r><w:rPr><w:lang w:val=\"pl-PL\"/></w:rPr><w:t>Baza danych:</w:t></w:r></w:p>'\n)\n$bazaPara.Range.InsertXML($bazaXml)\n\n# ---------------------------------------------------------------------\n# Step 3: rewrite the trailing \"Konwencja nazewnicza\" paragraph (merging\n# the split \"Konwe\"/\"ncja...\" runs and dropping the bookmark from here)\n# and append the new \"Cel aplikacji\" paragraphs after it, with the\n# _GoBack bookmark re-attached to the new final run.\n# ---------------------------------------------------------------------\n$n = $d.Paragraphs.Count\n$lastPara = $d.Paragraphs($n)\n$lastText = $lastPara.Range.Text\nif ($lastText -notlike \"*Konwe*\") {\n  throw \"Could not find trailing 'Konwencja nazewnicza' paragraph\"\n}\n\n$tailXml = New-PkgXml (\n  '<w:p><w:pPr><w:rPr><w:lang w:val=\"pl-PL\"/></w:rPr></w:pPr>' +\n  '<w:r><w:rPr><w:lang w:val=\"pl-PL\"/></w:rPr><w:t xml:space=\"preserve\">Konwencja nazewnicza zgodna ze standardami pisania kodu w </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:rPr><w:lang w:val=\"pl-PL\"/></w:rPr><w:t>Pythonie</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:rPr><w:lang w:val=\"pl-PL\"/></w:rPr><w:t xml:space=\"preserve\"> (PEP8).</w:t></w:r>' +\n  '</w:p>' +\n\n  '<w:p><w:pPr><w:rPr><w:lang w:val=\"pl-PL\"/></w:rPr></w:pPr>' +\n  '<w:r><w:rPr><w:lang w:val=\"pl-PL\"/></w:rPr><w:t xml:space=\"preserve\">Nie definiujemy konwencji  programowania dla bazy danych poniewa\u017c dla generowanych modeli w </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:rPr><w:lang w:val=\"pl-PL\"/></w:rPr><w:t>Django</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:rPr><w:lang w:val=\"pl-PL\"/></w:rPr><w:t xml:space=\"preserve\"> nie istnieje potrzeba pisania w j\u0119zyku SQL.</w:t></w:r>' +\n  '</w:p>' +\n\n  '<w:p><w:pPr><w:rPr><w:lang w:val=\"pl-PL\"/></w:rPr></w:pPr></w:p>' +\n\n  '<w:p><w:pPr><w:rPr><w:sz w:val=\"28\"/><w:lang w:val=\"pl-PL\"/></w:rPr></w:pPr>' +\n  '<w:r><w:rPr><w:sz w:val=\"28\"/><w:lang w:val=\"pl-PL\"/></w:rPr><w:t>5. Cel aplikacji</w:t></w:r>' +\n  '</w:p>' +\n\n  '<w:p><w:pPr><w:rPr><w:sz w:val=\"28\"/><w:lang w:val=\"pl-PL\"/></w:rPr></w:pPr></w:p>' +\n\n  '<w:p><w:pPr><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/><w:lang w:val=\"pl-PL\"/></w:rPr></w:pPr>' +\n  '<w:r><w:rPr><w:sz w:val=\"24\"/><w:lang w:val=\"pl-PL\"/></w:rPr>' +\n  '<w:t xml:space=\"preserve\">Aplikacja ma na celu informowa\u0107 potencjalnych zawodnik\u00f3w o mo\u017cliwo\u015bci przyst\u0105pienia do zawod\u00f3w. Wspomaga\u0107 rejestracj\u0119 uczestnik\u00f3w. Informowa\u0107 na bie\u017c\u0105co o tym jakie s\u0105 wyniki </w:t></w:r>' +\n  '<w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/><w:lang w:val=\"pl-PL\"/></w:rPr>' +\n  '<w:t xml:space=\"preserve\">zawod\u00f3w. Umo\u017cliwia\u0107 agregacj\u0119 tre\u015bci zwi\u0105zanych z zawodami. </w:t></w:r>' +\n  '</w:p>' +\n\n  '<w:p><w:pPr><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/><w:lang w:val=\"pl-PL\"/></w:rPr></w:pPr></w:p>' +\n\n  '<w:p><w:pPr><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/><w:lang w:val=\"pl-PL\"/></w:rPr></w:pPr></w:p>' +\n\n  '<w:p><w:pPr><w:rPr><w:sz w:val=\"28\"/><w:szCs w:val=\"24\"/><w:lang w:val=\"pl-PL\"/></w:rPr></w:pPr>' +\n  '<w:r><w:rPr><w:sz w:val=\"28\"/><w:szCs w:val=\"24\"/><w:lang w:val=\"pl-PL\"/></w:rPr><w:t xml:space=\"preserve\">6. </w:t></w:r>' +\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>' +\n  '</w:p>' +\n\n  '<w:p><w:pPr><w:rPr><w:lang w:val=\"pl-PL\"/></w:rPr></w:pPr></w:p>'\n)\n\n$lastPara.Range.InsertXML($tailXml)\n\n\"done\"\n"}
